$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the values in A22:A51 while keeping their existing cell style/formatting.
$ws.Range("A22:A51").ClearContents()

# Update the view: scroll so A36 is the top-left visible cell, and select A51.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A51").Select()
